# início das queries de análise serial
# Update category labels (column A) on the "max-arrecad" and "tx-sucesso"
# sheets to reflect the new serial-query tie-break ordering. The B column
# (values) is unchanged; only which category name sits on which row moves.

$wb = $excel.ActiveWorkbook

# --- Sheet: max-arrecad -----------------------------------------------
$wsMax = $wb.Worksheets.Item("max-arrecad")

$wsMax.Range("A2").Value = "humor"
$wsMax.Range("A3").Value = "folclore"
$wsMax.Range("A4").Value = "religiosidade"
$wsMax.Range("A5").Value = "terror"
$wsMax.Range("A6").Value = "fantasia"
$wsMax.Range("A7").Value = "questoes_genero"
$wsMax.Range("A8").Value = "fiq"
$wsMax.Range("A9").Value = "ficcao_cientifica"

$wsMax.Range("A13").Value = "hqmix"
$wsMax.Range("A14").Value = "angelo_agostini"

# --- Sheet: tx-sucesso --------------------------------------------------
$wsTx = $wb.Worksheets.Item("tx-sucesso")

$wsTx.Range("A2").Value = "angelo_agostini"
$wsTx.Range("A3").Value = "ccxp"
$wsTx.Range("A4").Value = "hqmix"
$wsTx.Range("A5").Value = "saloes_humor"

$wsTx.Range("A8").Value = "erotismo"
$wsTx.Range("A9").Value = "questoes_genero"

$wsTx.Range("A18").Value = "herois"
$wsTx.Range("A19").Value = "religiosidade"
